$wb = $excel.ActiveWorkbook

# --- Rename & reorder sheets ---
# "Plan1" holds the existing "LRC - Wavelet" dataset; it becomes the second tab.
# "Plan2" is currently empty; it becomes the first tab "LRC" with a new summary table.
$wb.Worksheets.Item("Plan2").Name = "LRC"
$wb.Worksheets.Item("Plan1").Name = "LRC - Wavelet"
$wb.Worksheets.Item("LRC").Move($wb.Worksheets.Item(1))

# Re-fetch the sheet handles by name now that positions have changed -
# stale handles obtained before a Move() can resolve to the wrong sheet.
$wsLrc = $wb.Worksheets.Item("LRC")
$wsWavelet = $wb.Worksheets.Item("LRC - Wavelet")

# --- Populate the new "LRC" sheet ---
# (16.44140625 char-units is the authored width; the nearest value this
# engine's column-width quantization can reach is 16.5.)
$wsLrc.Columns("A").ColumnWidth = 15.6

$wsLrc.Range("A1").Value = "LRC - 100 holdout experimento (50% treino - 50% teste)"

$wsLrc.Range("B2").Value = "ARFaces"
$wsLrc.Range("C2").Value = "YaleB"
$wsLrc.Range("D2").Value = "Gtech"
$wsLrc.Range("E2").Value = "ORL"
$wsLrc.Range("F2").Value = "EssexFaces"

$wsLrc.Range("A3").Value = "Taxa de Acerto"
$wsLrc.Range("B3").Value = 0.81005384615384601
$wsLrc.Range("C3").Value = 0.929095394736842
$wsLrc.Range("D3").Value = 0.77851428571428605
$wsLrc.Range("E3").Value = 0.9587
$wsLrc.Range("F3").Value = 0.92076388888888805

$wsLrc.Range("A4").Value = "Desvio Padrao"
$wsLrc.Range("B4").Value = 0.0121941415948623
$wsLrc.Range("C4").Value = 0.0077335851447940598
$wsLrc.Range("D4").Value = 0.019583290240613602
$wsLrc.Range("E4").Value = 0.016199825413699499
$wsLrc.Range("F4").Value = 0.0091967057889631998

# --- Update titles on the "LRC - Wavelet" sheet ---
$wsWavelet.Range("A1").Value = "LRC - Wavelet - 100 holdouts experimento (50% treino - 50% teste)"
$wsWavelet.Range("A22").Value = "GTech - Waveletfaces lvl3"

# --- Update selection / active sheet ---
$wsWavelet.Range("J8").Select()
$wsLrc.Activate()
$wsLrc.Range("A2").Select()
